$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPtFM")

# New row 7: electrolysis with guaranteed clean electricity (copies row 2)
$ws.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$ws.Range("B7").Formula = "=B2"
$ws.Range("C7:K7").Formula = "=C2"

# New row 8: natural gas reforming with CCS (copies row 3)
$ws.Range("A8").Value = "natural gas reforming with CCS"
$ws.Range("B8").Formula = "=B3"
$ws.Range("C8:K8").Formula = "=C3"

# Sheet "About": A7 previously had an explicit (but visually default) style
# applied; clearing the bold flag removes the now-redundant style so the
# cell reverts to the default formatting.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A7").Font.Bold = $false
